# Add a new "forecast as of 2020-05-14" column (U) and a new row (33, for
# 2020-05-14 itself) to both the "cases" and "deaths" sheets, mirroring the
# existing layout:
#   - row 1 holds target-date headers (shared strings) in B1:T1 -> now U1
#   - column A holds "as of" row-date labels A2:A32 -> now also A33
#   - data cells are plain numbers (or blank)

$wb = $excel.ActiveWorkbook

# A range with the sheet's default (unstyled) cell - used below to strip any
# incidental number-format styling back off cells after we coerce a
# date-shaped string into literal text.
function Get-DefaultStyle($ws) {
    return $ws.Range("A1").Style
}

# Force a literal (non date-auto-converted) text value into a cell: Excel's
# COM layer infers dates from strings like "2020-05-14" unless the cell is
# already formatted as Text; we flip it to Text, assign, then restore the
# worksheet's default style so no stray per-cell formatting is left behind.
function Set-TextValue($ws, $cellAddr, $text, $defaultStyle) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $defaultStyle
}

# Touch a cell so an explicit (empty) <c> node is persisted for it, matching
# rows that already enumerate every column even when blank. Re-applying the
# (non-bold) font is a no-op visually/structurally but forces the engine to
# materialize the cell.
function Set-BlankCell($ws, $cellAddr) {
    $ws.Range($cellAddr).Font.Bold = $false
}

$sheetNames = @("cases", "deaths")

# New U-column (forecast horizon "2020-05-14") values per data row, and the
# new trailing row 33 ("as of 2020-05-14") U-value + B19 fill-in, per sheet.
$caseData = @{
    B19 = 9453
    U20 = 10314
    U21 = 10570
    U22 = 11993
    U23 = 12236
    U24 = 13402
    U25 = 14574
    U26 = 15657
    U27 = 16608
    U28 = 17989
    U29 = 19559
    U30 = 21118
    U31 = 22674
    U32 = 23770
    U33 = 24774
}

$deathData = @{
    B19 = 854
    U20 = 915
    U21 = 940
    U22 = 1065
    U23 = 1088
    U24 = 1191
    U25 = 1295
    U26 = 1390
    U27 = 1476
    U28 = 1597
    U29 = 1738
    U30 = 1874
    U31 = 2010
    U32 = 2107
    U33 = 2197
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $defaultStyle = Get-DefaultStyle $ws

    if ($sheetName -eq "cases") {
        $data = $caseData
    } else {
        $data = $deathData
    }

    # --- Header row: new column U1, reusing the existing "2020-04-30" label
    # (shared string index 20) just like the source workbook does. ---
    Set-TextValue $ws "U1" "2020-04-30" $defaultStyle

    # --- Existing rows 2-18: add an empty U cell to round out the row. ---
    for ($r = 2; $r -le 18; $r++) {
        Set-BlankCell $ws ("U" + $r)
    }

    # --- Row 19: B19 was blank, now has a real forecast value; also add the
    # trailing empty U19 cell. ---
    $ws.Range("B19").Value = $data["B19"]
    Set-BlankCell $ws "U19"

    # --- Rows 20-32: append the new U-column forecast values. ---
    for ($r = 20; $r -le 32; $r++) {
        $ws.Range("U" + $r).Value = $data["U" + $r]
    }

    # --- New row 33 ("as of 2020-05-14"): label in A33 (new shared string),
    # blank B33:T33, and the U33 forecast value. ---
    Set-TextValue $ws "A33" "2020-05-14" $defaultStyle
    for ($c = 2; $c -le 20; $c++) {
        $colLetter = [char](64 + $c)
        Set-BlankCell $ws ($colLetter + "33")
    }
    $ws.Range("U33").Value = $data["U33"]
}
